$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("120:122").Delete()
